$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Shift the "Lower Right Cell" references for the scenario index blocks
# from row 43 to row 46 (D5:D11), to reflect extra scenario rows added.
$ws.Range("D5").Value = "A46"
$ws.Range("D6").Value = "B46"
$ws.Range("D7").Value = "C46"
$ws.Range("D8").Value = "G46"
$ws.Range("D9").Value = "H46"
$ws.Range("D10").Value = "I46"
$ws.Range("D11").Value = "J46"

# Update the selection shown in the saved view.
$ws.Range("D5:D11").Select()
